$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.879.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.715.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +6.94%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "421.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.703.06"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +6.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.647"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.767"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.188"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +19.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000402"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +76.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.299.25"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.81%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.57"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.703.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.86%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.880.76"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "451.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.28"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +15.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.80"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.33"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.13"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.13"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.90%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0487"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0749"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +18.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "147.92"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +16.57%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.68"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.29%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.07"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +22.01%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.30"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.28%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.165"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +19.01%  "
